# Adds the final set of quests for the new event map:
# "Infiltrating the Enemy" -> "A word from our prophet" -> "The judges of time" -> "Celestial Hunt"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 181 (id 180)
$ws.Cells.Item(181, 1).Value = 180
$ws.Cells.Item(181, 2).Value = 'Infiltrating the Enemy'
$ws.Cells.Item(181, 3).Value = 'FederationChancellor'
$ws.Cells.Item(181, 5).Value = 'A man named: Fliniguss'
$ws.Cells.Item(181, 6).Value = 'The Child Returns'
$ws.Cells.Item(181, 7).Value = 20
$ws.Cells.Item(181, 8).Value = 'Twisted Tree Branch'
$ws.Cells.Item(181, 9).Value = 550000
$ws.Cells.Item(181, 10).Value = 550000
$ws.Cells.Item(181, 11).Value = 150000000000
$ws.Cells.Item(181, 12).Value = 700000
$ws.Cells.Item(181, 14).Value = 500000
$ws.Cells.Item(181, 15).Value = 500000
$ws.Cells.Item(181, 16).Value = 2000000000
$ws.Cells.Item(181, 17).Value = 175000
$ws.Cells.Item(181, 21).Value = 'Looking for the Jester'
$ws.Cells.Item(181, 27).Value = 'You venture to Federation City in the far north to rendezvous with a man known only as the Federation Chancellor. Disguised under the pretense of investigating corruption within the upper echelons of The Federation, possibly implicating—as revealed by The Jester of Time, Fliniguss'' absence—his potential involvement.<br /> <br /> The Federation Chancellor acquiesces to your meeting, accompanied by representatives from the church. You''ve heard rumors of the twisted and sinister deeds attributed to both The Federation and The Church: the <br /> carnage, death, enslavement, and even more abhorrent horrors lurking within one''s psyche, forbidden from surfacing.<br /> <br /> This city appears more technologically advanced than what you''re accustomed to, contrasting the stone and wood constructions of other cities. Here, metal, steam, and gears prevail, indicative of a higher level of <br /> advancement. Giant airships dominate the skies, while individuals rush about the streets in modern attire.<br /> <br /> You feel conspicuously out of place, not just in this era but also in this locale.<br /> <br /> Navigating through one of the bustling main streets, teeming with small pubs, eateries, and shops brimming with peculiar curiosities, you pause to peer into a shop window. Suddenly, a voice startles you from behind: <br /> "Steampunk is the design, or so they say. Contraptions fueled by steam are presumably crafted by mischievous individuals. Whatever the case, these toys captivate my children for hours."<br /> <br /> You turn to behold a sharply dressed man.<br /> <br /> "I am the Federation Chancellor you seek. My office lies across the street. A bishop and clergy await us. Shall we proceed?" he inquires.<br /> <br /> You nod in affirmation and trail behind him as he strides across the street and into the building.'
$ws.Cells.Item(181, 28).Value = 'You find yourself seated at a substantial oak table alongside a Bishop of The Church, a clergy member of The Church, and the Federation Chancellor. Their exchange is fraught with argument and contention:<br /> <br /> "It''s preposterous! These are sheer fabrications!" the Bishop bellows across the table.<br /> <br /> "I''m unfamiliar with you. Whom do you represent?" the clergy member directs towards you after a pause.<br /> <br /> You gaze at both men, fabricating a story with practiced ease. You claim allegiance to Fliniguss and assert that his recent disappearance prompted your inquiry into matters that hinted at corruption. Your search led you <br /> to unsettling revelations, culminating in the discovery of his lifeless body.<br /> <br /> "Nonsense! I refuse to entertain such baseless claims. If Fliniguss were missing or worse, the Arch Bishop would have been informed by now," the Bishop retorts vehemently.<br /> <br /> "I propose we subject this heretic to the judgment of the prophet. He shall discern truth from falsehood," the Bishop declares, slamming his hand on the table. "No! We must deliver him to The Church, allowing divine scrutiny to reveal the purity of his soul," counters another.<br /> <br /> "I concur. Why entrust him to a prophet the church has disavowed?" A query begins to form in your mind, but the Chancellor interjects.<br /> <br /> "Fliniguss would choose this course. I understand his mindset. He would opt for the prophet in the forsaken village, believing wholeheartedly in the man''s divine connection," the Chancellor asserts.<br /> <br /> "He''s a heretic, a deceiver, and it''s common knowledge that the Arch Bishop turns a blind eye to Fliniguss and his association with this individual," the clergy member interjects.<br /> <br /> A heavy silence descends upon the table, but you declare your willingness to proceed with this plan. The Chancellor''s reasoning resonates; Fliniguss would indeed opt for such a path. Let us, if nothing else, honor that <br /> aspect of the man''s memory.<br /> <br /> "A shrewd maneuver indeed. We shall see who emerges victorious," the Chancellor remarks with a smile. The others recline in their chairs, scrutinizing you intently.<br /> <br /> A game is set in motion, one that promises either revelation or demise.'
$ws.Cells.Item(181, 30).Value = 'A man named: Fliniguss'
$ws.Cells.Item(181, 31).Value = 8
$ws.Cells.Item(181, 32).Value = 'FederationChancellor'
$ws.Cells.Item(181, 33).Value = 5

# Row 182 (id 181)
$ws.Cells.Item(182, 1).Value = 181
$ws.Cells.Item(182, 2).Value = 'A word from our prophet'
$ws.Cells.Item(182, 3).Value = 'ProphetofTheChurch'
$ws.Cells.Item(182, 12).Value = 150000
$ws.Cells.Item(182, 14).Value = 500000
$ws.Cells.Item(182, 15).Value = 500000
$ws.Cells.Item(182, 16).Value = 2000000000
$ws.Cells.Item(182, 17).Value = 175000
$ws.Cells.Item(182, 21).Value = 'Infiltrating the Enemy'
$ws.Cells.Item(182, 23).Value = 'Twisted Memories'
$ws.Cells.Item(182, 24).Value = 5
$ws.Cells.Item(182, 27).Value = 'After days of travel, you stumble upon an abandoned chapel nestled in the heart of a desolate village. Not a soul stirs; the pervasive silence weighs heavily, suffocating in its loneliness.<br /> <br /> Venturing inside the dilapidated chapel, its structure half-burnt and mostly crumbled, you encounter a solitary figure seated upon a remaining pew. As he rises to greet you, his appearance strikes a balance between age and vitality – weathered by time yet possessing a serene countenance.<br /> <br /> "Child, I understand your purpose here," he intones softly.<br /> <br /> You convey your mission: persuading the Chancellor and the Bishop of Fliniguss''s demise, unable to disclose further details but hoping to glean insights into their next course of action now that their organization''s leader is no more.<br /> <br /> "He cannot have perished, child. That does not align with God''s plan for him," he responds with gentle certainty. Intrigued, you inquire how he is so certain.<br /> <br /> "Silly child. God speaks to me. He also reveals a trial awaiting you – four judges who roam the corridors of the old church. Judges tasked with testing your mettle in confronting harsh truths," he explains, a faint chuckle <br /> escaping his lips as he steps past you and exits the chapel into the desolate village.<br /> <br /> "This place," he begins, gesturing around him, "once teemed with life. Now, it is but a testament to death. Each inhabitant faced judgment, yet few prevailed."<br /> <br /> He turns to regard you intently.<br /> <br /> "But he survived. He is not among the deceased," he asserts.'
$ws.Cells.Item(182, 28).Value = 'You find yourself spending a considerable amount of time within the village, engaging in activities such as hunting local wildlife and assisting the prophet with his enigmatic "chores."<br /> <br /> Despite your efforts, answers remain elusive – a man ostensibly deceased yet alive, the enigmatic Jester of Time, and elusive dreams of The Creator.<br /> <br /> Your attempts at espionage against The Federation are stymied by a man who claims divine communication. The concept of God perplexes you – is it akin to the gods and celestial beings of this realm? Or is it yet <br /> another unseen force pulling strings behind the scenes?<br /> <br /> Questions multiply, answers remain elusive. This individual leads you in circles, offering cryptic hints but no clear solutions. Yet, the prospect of facing the judges intrigues you – perhaps therein lies a semblance of clarity amid the chaos.'
$ws.Cells.Item(182, 31).Value = 8
$ws.Cells.Item(182, 32).Value = 'ProphetofTheChurch'
$ws.Cells.Item(182, 33).Value = 5

# Row 183 (id 182)
$ws.Cells.Item(183, 1).Value = 182
$ws.Cells.Item(183, 2).Value = 'The judges of time'
$ws.Cells.Item(183, 3).Value = 'ProphetofTheChurch'
$ws.Cells.Item(183, 14).Value = 500000
$ws.Cells.Item(183, 15).Value = 500000
$ws.Cells.Item(183, 16).Value = 2000000000
$ws.Cells.Item(183, 17).Value = 175000
$ws.Cells.Item(183, 21).Value = 'A word from our prophet'
$ws.Cells.Item(183, 27).Value = 'Inquiring with the prophet about the judges yields cryptic insight.<br /> <br /> "There are four of them, child. They judge the essence of your being, not through brute force like your hunts here, but by surviving through the guidance of God''s voice and grace. It is through divine guidance that you <br /> shall find the answers you seek," he explains.<br /> <br /> You scoff, claiming he knows nothing of your inquiries.<br /> <br /> "I am aware of your pursuit of the prince, your encounters with the jester, and your struggles with the illusions of this world. But what you fail to see is the bigger picture – you are as blind as The Child, as The Creator, <br /> even as the Madam who believes she controls memories lost to the twisted whims of time," he asserts, fixing his gaze upon you.<br /> <br /> After a pregnant pause, he continues, "I perceive a yearning within you, a quest for something beyond your comprehension. And I know you are an inadequate spy. Your deceit brought you this far, but The Federation <br /> holds no answers for you."<br /> <br /> Not even regarding Fliniguss''s true demise? Or the thwarted ascension? You are aware that this is a distorted memory from a deluded era.<br /> <br /> "You are not of this time, child. You possess knowledge beyond your era and recall memories obscured by fog, yet clarified by your words," he remarks before turning away.<br /> <br /> As you call out to him, he continues walking, his destination shrouded in mystery. Where is he headed?'
$ws.Cells.Item(183, 28).Value = 'You wait in anticipation, but days pass without the prophet''s presence. The rhythm of rain and silence envelops the land.<br /> <br /> One day, as you return to the chapel after foraging for herbs and game, you spot him standing at the entrance. Observing him for a moment, you approach and inquire about his absence.<br /> <br /> "I sought revelation, child. There is much to discuss," he replies, his hand resting on your shoulder.<br /> <br /> Though he claims not to recall you, he acknowledges awareness of your kind – transient beings who emerge from obscurity, probing with countless inquiries before fading away once more. According to his divine insight, <br /> you possess the capacity to confront the judges of time – arbiters of your essence, body, and psyche; arbiters of your anguish.<br /> <br /> Eagerly, you inquire how to face these judges, hoping for further elucidation.<br /> <br /> "They hold no answers," he responds with a smile.'
$ws.Cells.Item(183, 31).Value = 8
$ws.Cells.Item(183, 32).Value = 'ProphetofTheChurch'
$ws.Cells.Item(183, 33).Value = 10

# Row 184 (id 183)
$ws.Cells.Item(184, 1).Value = 183
$ws.Cells.Item(184, 2).Value = 'Celestial Hunt'
$ws.Cells.Item(184, 3).Value = 'ProphetofTheChurch'
$ws.Cells.Item(184, 14).Value = 500000
$ws.Cells.Item(184, 15).Value = 500000
$ws.Cells.Item(184, 16).Value = 2000000000
$ws.Cells.Item(184, 17).Value = 175000
$ws.Cells.Item(184, 21).Value = 'The judges of time'
$ws.Cells.Item(184, 27).Value = 'Early one morning, as you approach the ruined chapel to meet the prophet, he greets you with news delivered by carrier pigeon.<br /> <br /> "They inquired about your progress," he reveals. "Whether they should dispatch the guard. I responded by informing them that I have sent you to the Judges."<br /> <br /> You express frustration, recounting the Bishop and the clergy''s dismissal of your credibility, branding you a heretic.<br /> <br /> "They hurl insults, they cast their stones. Yet God speaks to me, and the Arch Bishop, he discerns the truth. Despite murmurs of doubt, I perceive the state of his soul. It has undergone judgment," he declares, settling on <br /> the chapel steps.<br /> <br /> Seizing the opportunity, you inquire about the whereabouts of the judges and how to reach them, how to engage in this confrontation.<br /> <br /> "First, child, you must survive," he cryptically replies.<br /> <br /> Survive? Survive what?'
$ws.Cells.Item(184, 28).Value = '"There exist creatures in this realm more fearsome than the quarry of your hunts for bones and teeth," the prophet remarks solemnly.<br /> <br /> You inform him of these beings, known as celestials. Their origins remain a mystery, yet you''ve heard of individuals possessing the magic to summon them. Moreover, you recount instances when they emerge from their <br /> hidden lairs, fueled by rage, hatred, and the depths of darkness they harbor.<br /> <br /> "Abominable creations. They must be vanquished! We must restore peace to God''s earth," he declares with fervor.<br /> <br /> However, you cannot help but ponder that true peace seems an elusive dream as long as darkness persists.'
$ws.Cells.Item(184, 31).Value = 8
$ws.Cells.Item(184, 32).Value = 'ProphetofTheChurch'
$ws.Cells.Item(184, 33).Value = 15

# gold_cost column widened to fit the new 150,000,000,000 value (12 digits)
$ws.Columns.Item(11).AutoFit()
